# Automatische test-sync: 2025-08-28 18:16:50
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 9)
$ws.Range("A9").Value = "Retour status"
$ws.Range("B9").Value = "mailmind.test@zohomail.eu"
$ws.Range("D9").Value = "Retour / Terugbetaling"
$ws.Range("F9").Value = "2025-08-28 18:16:14"
$ws.Range("G9").Value = "Nee"
$ws.Range("H9").Value = "Ja"
$ws.Range("I9").Value = "Nee"
$ws.Range("J9").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range("$col" + "2:" + "$col" + "8")
    $newRange = $ws.Range("$col" + "2:" + "$col" + "9")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 8
